$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 162, shifting existing rows 162-179 down to 163-180
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new weekly record
$ws.Range("A162").Value = 3
$ws.Range("B162").Value = "Femacal de La Calera"
$ws.Range("C162").Value = "Coquimbo"
$ws.Range("D162").Value = 44449
$ws.Range("E162").Value = 5
$ws.Range("F162").Value = 100112012
$ws.Range("G162").Value = "Espinaca"
$ws.Range("H162").Value = "Sin especificar"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 230
$ws.Range("K162").Value = 2500
$ws.Range("L162").Value = 2800
$ws.Range("M162").Value = 2643
$ws.Range("N162").Value = '$/docena de atados (3 kilos)'
$ws.Range("O162").Value = "Provincia de Quillota"
$ws.Range("P162").Value = 881
$ws.Range("Q162").Value = 3
$ws.Range("R162").Value = "Hortaliza"
